# The commit swaps the presentation's theme away from the custom
# "Integral" design and onto the stock PowerPoint "Office Theme"
# (Design tab -> Themes -> Office). In the saved OOXML this shows up as
# the clrScheme / fontScheme / fmtScheme content of ppt/theme/theme1.xml
# (the theme used by the one slide master) becoming the default Office
# theme definition, while the old Integral definition is preserved in
# ppt/theme/theme2.xml (used by the notes master).
#
# Apply the new "Office" theme color palette to the slide master's
# theme color scheme. Colors are expressed as decimal BGR integers
# (PowerPoint's RGB long format: R + G*256 + B*65536) since this host
# does not expose a RGB() helper.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme palette (standard Office 2013+ default theme)
$tcs.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # Dark 1   - 000000
$tcs.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # Light 1  - FFFFFF
$tcs.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # Dark 2   - 44546A
$tcs.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # Light 2  - E7E6E6
$tcs.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # Accent 1 - 5B9BD5
$tcs.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # Accent 2 - ED7D31
$tcs.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # Accent 3 - A5A5A5
$tcs.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # Accent 4 - FFC000
$tcs.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # Accent 5 - 4472C4
$tcs.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # Accent 6 - 70AD47
$tcs.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # Hyperlink          - 0563C1
$tcs.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # Followed Hyperlink - 954F72
